$wb = $excel.ActiveWorkbook

# ALC!row 98
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1312.8125
$ws.Range("I98").Value = 1042.0834
$ws.Range("J98").Value = 2125
$ws.Range("K98").Value = 1042.0834
$ws.Range("L98").Value = 2125
$ws.Range("M98").Value = 455.9166
$ws.Range("N98").Value = -5121

# ALC!row 112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 58825652
$ws.Range("J112").Value = 2462.8572
$ws.Range("L112").Value = 7388.571599999999
$ws.Range("N112").Value = -9604.571599999999

# ALC!row 122
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 1312.8125
$ws.Range("I122").Value = 1042.0834
$ws.Range("J122").Value = 2125
$ws.Range("K122").Value = 3126.2502
$ws.Range("L122").Value = 6375
$ws.Range("M122").Value = -676.2501999999999
$ws.Range("N122").Value = -11275

# ALC!row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 614301.3
$ws.Range("I132").Value = 1681.9508
$ws.Range("J132").Value = 2581131.8
$ws.Range("K132").Value = 5045.8524
$ws.Range("L132").Value = 7743395.399999999
$ws.Range("M132").Value = -2515.8524
$ws.Range("N132").Value = -7748455.399999999

# ALC!row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1539804.9
$ws.Range("I137").Value = 1852726.4
$ws.Range("J137").Value = 3645.182
$ws.Range("K137").Value = 5558179.199999999
$ws.Range("L137").Value = 10935.546
$ws.Range("M137").Value = -5555629.199999999
$ws.Range("N137").Value = -16035.546

# ALC!row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2010084.1
$ws.Range("I138").Value = 1178.8392
$ws.Range("J138").Value = 6176702.5
$ws.Range("K138").Value = 3536.5176
$ws.Range("L138").Value = 18530107.5
$ws.Range("M138").Value = 1603.4824
$ws.Range("N138").Value = -18540387.5

# ALC!row 141
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4418.12
$ws.Range("I141").Value = 2836
$ws.Range("J141").Value = 8486.429
$ws.Range("K141").Value = 8508
$ws.Range("L141").Value = 25459.287
$ws.Range("M141").Value = -3328
$ws.Range("N141").Value = -35819.287

# ARM!row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18783.715
$ws.Range("I32").Value = 16825.402
$ws.Range("J32").Value = 25964.191
$ws.Range("K32").Value = 16825.402
$ws.Range("L32").Value = 25964.191
$ws.Range("M32").Value = -16538.402
$ws.Range("N32").Value = -26538.191

# ARM!row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1166.5
$ws.Range("I45").Value = 1022.8333
$ws.Range("J45").Value = 1597.5
$ws.Range("K45").Value = 1022.8333
$ws.Range("L45").Value = 1597.5
$ws.Range("M45").Value = -645.8333
$ws.Range("N45").Value = -2351.5

# ARM!row 61
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 15182957
$ws.Range("I61").Value = 18887794
$ws.Range("J61").Value = 78622.16
$ws.Range("K61").Value = 18887794
$ws.Range("L61").Value = 78622.16
$ws.Range("M61").Value = -18887582
$ws.Range("N61").Value = -79046.16

# ARM!row 74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 9864257
$ws.Range("I74").Value = 13931900
$ws.Range("J74").Value = 101914.7
$ws.Range("K74").Value = 13931900
$ws.Range("L74").Value = 101914.7
$ws.Range("M74").Value = -13931026
$ws.Range("N74").Value = -103662.7

# ARM!row 77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 9864257
$ws.Range("I77").Value = 13931900
$ws.Range("J77").Value = 101914.7
$ws.Range("K77").Value = 69659500
$ws.Range("L77").Value = 509573.5
$ws.Range("M77").Value = -69655132
$ws.Range("N77").Value = -518309.5

# ARM!row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 44533.895
$ws.Range("I132").Value = 31978.094
$ws.Range("J132").Value = 69645.5
$ws.Range("K132").Value = 95934.28200000001
$ws.Range("L132").Value = 208936.5
$ws.Range("M132").Value = -93404.28200000001
$ws.Range("N132").Value = -213996.5

# ARM!row 136
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 15182957
$ws.Range("I136").Value = 18887794
$ws.Range("J136").Value = 78622.16
$ws.Range("K136").Value = 56663382
$ws.Range("L136").Value = 235866.48
$ws.Range("M136").Value = -56660832
$ws.Range("N136").Value = -240966.48

# BSM!row 92
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 10400
$ws.Range("J92").Value = 10400
$ws.Range("L92").Value = 10400
$ws.Range("N92").Value = -15392

# BSM!row 95
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H95").Value = 15624
$ws.Range("J95").Value = 15624
$ws.Range("L95").Value = 15624
$ws.Range("N95").Value = -21116

# BSM!row 100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 20000
$ws.Range("J100").Value = 20000
$ws.Range("L100").Value = 20000
$ws.Range("N100").Value = -22164

# BSM!row 134
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3714.3
$ws.Range("I134").Value = 2029.1
$ws.Range("J134").Value = 5399.5
$ws.Range("K134").Value = 6087.299999999999
$ws.Range("L134").Value = 16198.5
$ws.Range("M134").Value = -3552.299999999999
$ws.Range("N134").Value = -21268.5

# CRP!row 31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1410.8833
$ws.Range("I31").Value = 787.5282999999999
$ws.Range("K31").Value = 787.5282999999999
$ws.Range("M31").Value = -492.5282999999999

# CRP!row 34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1410.8833
$ws.Range("I34").Value = 787.5282999999999
$ws.Range("K34").Value = 787.5282999999999
$ws.Range("M34").Value = -585.5282999999999

# CRP!row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 18462.508
$ws.Range("I132").Value = 1367.6744
$ws.Range("J132").Value = 64404.875
$ws.Range("K132").Value = 4103.023200000001
$ws.Range("L132").Value = 193214.625
$ws.Range("M132").Value = -1573.023200000001
$ws.Range("N132").Value = -198274.625

# CRP!row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 15901.768
$ws.Range("I134").Value = 1049.4237
$ws.Range("J134").Value = 78493.78999999999
$ws.Range("K134").Value = 3148.2711
$ws.Range("L134").Value = 235481.37
$ws.Range("M134").Value = -613.2710999999999
$ws.Range("N134").Value = -240551.37

# CUL!row 103
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 3666.6428
$ws.Range("I103").Value = 461
$ws.Range("J103").Value = 5447.5557
$ws.Range("K103").Value = 1383
$ws.Range("L103").Value = 16342.6671
$ws.Range("M103").Value = -504
$ws.Range("N103").Value = -18100.6671

# CUL!row 123
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H123").Value = 2105.8572
$ws.Range("I123").Value = 1735
$ws.Range("J123").Value = 3033
$ws.Range("K123").Value = 5205
$ws.Range("L123").Value = 9099
$ws.Range("M123").Value = -2755
$ws.Range("N123").Value = -13999

# CUL!row 129
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2875477.2
$ws.Range("I129").Value = 1560.8182
$ws.Range("J129").Value = 4631759.5
$ws.Range("K129").Value = 4682.4546
$ws.Range("L129").Value = 13895278.5
$ws.Range("M129").Value = 317.5454
$ws.Range("N129").Value = -13905278.5

# CUL!row 130
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 2903
$ws.Range("I130").Value = 900
$ws.Range("J130").Value = 3125.5557
$ws.Range("K130").Value = 2700
$ws.Range("L130").Value = 9376.667099999999
$ws.Range("M130").Value = 2320
$ws.Range("N130").Value = -19416.6671

# CUL!row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1488.3903
$ws.Range("J131").Value = 1642.3334
$ws.Range("L131").Value = 4927.0002
$ws.Range("N131").Value = -15007.0002

# CUL!row 133
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H133").Value = 6358
$ws.Range("I133").Value = 3695
$ws.Range("J133").Value = 7996.769
$ws.Range("K133").Value = 11085
$ws.Range("L133").Value = 23990.307
$ws.Range("M133").Value = -6025
$ws.Range("N133").Value = -34110.307

# CUL!row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 6655.4346
$ws.Range("I134").Value = 2497
$ws.Range("K134").Value = 7491
$ws.Range("M134").Value = -2421

# CUL!row 136
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 3158.3333
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 3158.3333
$ws.Range("K136").Value = 0
$ws.Range("L136").ClearContents()
$ws.Range("M136").Value = 9474.999899999999
$ws.Range("N136").Value = -19674.9999

# CUL!row 137
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 23467.592
$ws.Range("I137").Value = 922.8570999999999
$ws.Range("J137").Value = 31358.25
$ws.Range("K137").Value = 2768.5713
$ws.Range("L137").Value = 94074.75
$ws.Range("M137").Value = 2331.4287
$ws.Range("N137").Value = -104274.75

# CUL!row 139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 4072.1694
$ws.Range("I139").Value = 2032.9
$ws.Range("J139").Value = 5117.9487
$ws.Range("K139").Value = 6098.700000000001
$ws.Range("L139").Value = 15353.8461
$ws.Range("M139").Value = -958.7000000000007
$ws.Range("N139").Value = -25633.8461

# CUL!row 140
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2673.239
$ws.Range("I140").Value = 2782.2856
$ws.Range("K140").Value = 8346.856800000001
$ws.Range("M140").Value = -3166.856800000001

# CUL!row 141
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 10298.556
$ws.Range("I141").Value = 10298.556
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 30895.668
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -25715.668

# GSM!row 98
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 28762
$ws.Range("J98").Value = 28762
$ws.Range("L98").Value = 28762
$ws.Range("N98").Value = -34752

# GSM!row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 40028.75
$ws.Range("I132").Value = 29030.473
$ws.Range("J132").Value = 64774.875
$ws.Range("K132").Value = 87091.41900000001
$ws.Range("L132").Value = 194324.625
$ws.Range("M132").Value = -84561.41900000001
$ws.Range("N132").Value = -199384.625

# GSM!row 136
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 14195.723
$ws.Range("J136").Value = 14195.723
$ws.Range("L136").Value = 42587.169
$ws.Range("N136").Value = -47687.169

# GSM!row 138
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 62000
$ws.Range("J138").Value = 62000
$ws.Range("L138").Value = 62000
$ws.Range("N138").Value = -72280

# LTW!row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5114
$ws.Range("I7").Value = 1766.6666
$ws.Range("K7").Value = 1766.6666
$ws.Range("M7").Value = -1654.6666

# LTW!row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5114
$ws.Range("I126").Value = 1766.6666
$ws.Range("K126").Value = 5299.9998
$ws.Range("M126").Value = -2829.9998

# LTW!row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 19611.436
$ws.Range("I132").Value = 1021.7
$ws.Range("J132").Value = 69184.07000000001
$ws.Range("K132").Value = 3065.1
$ws.Range("L132").Value = 207552.21
$ws.Range("M132").Value = -535.1000000000004
$ws.Range("N132").Value = -212612.21

# LTW!row 136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 48771.26
$ws.Range("I136").Value = 27913.324
$ws.Range("K136").Value = 83739.97200000001
$ws.Range("M136").Value = -81189.97200000001

# WVR!row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 28265.479
$ws.Range("I132").Value = 18451.822
$ws.Range("J132").Value = 60592.824
$ws.Range("K132").Value = 55355.466
$ws.Range("L132").Value = 181778.472
$ws.Range("M132").Value = -52825.466
$ws.Range("N132").Value = -186838.472

# WVR!row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 52228.55
$ws.Range("I136").Value = 40749.28
$ws.Range("J136").Value = 71360.664
$ws.Range("K136").Value = 122247.84
$ws.Range("L136").Value = 214081.992
$ws.Range("M136").Value = -119697.84
$ws.Range("N136").Value = -219181.992
